# Fix a typo on the "Assignment Statement" slide: the comment on the
# "name := "Caleb";" example line said "length = 4" but "Caleb" has 5
# characters, so the comment should read "length = 5".

$p = $ppt.ActivePresentation

# Locate the slide / shape / paragraph that contains the typo instead of
# hard-coding indices, so the script is resilient to minor reordering.
$targetShape = $null
$targetParagraph = $null

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if (-not $shape.HasTextFrame) {
            continue
        }
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like '*name := "Caleb";*length = 4*') {
            $paraCount = $tr.Paragraphs().Count
            for ($paraIdx = 1; $paraIdx -le $paraCount; $paraIdx++) {
                $para = $tr.Paragraphs($paraIdx, 1)
                if ($para.Text -like '*name := "Caleb";*length = 4*') {
                    $targetShape = $shape
                    $targetParagraph = $para
                    break
                }
            }
        }
        if ($targetShape -ne $null) {
            break
        }
    }
    if ($targetShape -ne $null) {
        break
    }
}

$tr = $targetShape.TextFrame.TextRange
$para = $targetParagraph

# Replace just the "= 4" portion with "= 5" -- this mirrors the author
# selecting the characters "= 4" and retyping "= 5". PowerPoint splits
# the run at the edited characters while leaving the rest of the
# line's formatting/run untouched.
$needle = "= 4"
$offsetInParagraph = $para.Text.IndexOf($needle)
$absoluteStart = $para.Start + $offsetInParagraph

$editRange = $tr.Characters($absoluteStart, $needle.Length)
$editRange.Text = "= 5"

Write-Output ("Corrected line: " + $targetShape.TextFrame.TextRange.Characters($para.Start, $needle.Length + $offsetInParagraph).Text)
